# Add a new "2020" column (Q) to the Employment-in-informal-sector table,
# mirroring the formatting already used for the 2019 column (P), and move
# the active selection to N17 (matches the author's last click before save).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header cell: Q4 = 2020, formatted like the other year headers ---
$ws.Range("Q4").Value = 2020
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)

# --- Data rows: copy number formatting from the matching P cell, then set the value ---
$ws.Range("Q5").Value = 1586.9
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)

$ws.Range("Q7").Value = 1032.4000000000001
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)

$ws.Range("Q8").Value = 554.5
$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)

$ws.Range("Q10").Value = 77.099999999999994
$ws.Range("P10").Copy()
$ws.Range("Q10").PasteSpecial(-4122)

$ws.Range("Q11").Value = 460.7
$ws.Range("P11").Copy()
$ws.Range("Q11").PasteSpecial(-4122)

$ws.Range("Q12").Value = 466.6
$ws.Range("P12").Copy()
$ws.Range("Q12").PasteSpecial(-4122)

$ws.Range("Q13").Value = 316.7
$ws.Range("P13").Copy()
$ws.Range("Q13").PasteSpecial(-4122)

$ws.Range("Q14").Value = 203.6
$ws.Range("P14").Copy()
$ws.Range("Q14").PasteSpecial(-4122)

$ws.Range("Q15").Value = 57.8
$ws.Range("P15").Copy()
$ws.Range("Q15").PasteSpecial(-4122)

$ws.Range("Q16").Value = 4.3
$ws.Range("P16").Copy()
$ws.Range("Q16").PasteSpecial(-4122)

$ws.Range("Q18").Value = 542.79999999999995
$ws.Range("P18").Copy()
$ws.Range("Q18").PasteSpecial(-4122)

$ws.Range("Q19").Value = 1044.0999999999999
$ws.Range("P19").Copy()
$ws.Range("Q19").PasteSpecial(-4122)

$ws.Range("Q21").Value = 339.3
$ws.Range("P21").Copy()
$ws.Range("Q21").PasteSpecial(-4122)

$ws.Range("Q22").Value = 230.4
$ws.Range("P22").Copy()
$ws.Range("Q22").PasteSpecial(-4122)

$ws.Range("Q23").Value = 270.5
$ws.Range("P23").Copy()
$ws.Range("Q23").PasteSpecial(-4122)

$ws.Range("Q24").Value = 746.7
$ws.Range("P24").Copy()
$ws.Range("Q24").PasteSpecial(-4122)

# --- Blank separator-row cells (still get touched with formatting) ---
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)

$ws.Range("P9").Copy()
$ws.Range("Q9").PasteSpecial(-4122)

$ws.Range("P17").Copy()
$ws.Range("Q17").PasteSpecial(-4122)

$ws.Range("P20").Copy()
$ws.Range("Q20").PasteSpecial(-4122)

# --- Keep row 36's spacer row inside the used range (dimension ref) ---
$ws.Range("A36").Value2 = $ws.Range("A36").Value2

# --- Match the saved selection location ---
$ws.Range("N17").Select()
